# Auto-generated edit script for cryptos.xlsx price/volume refresh
# Updates Price (D) and Volume(1h) (E) columns for the latest snapshot,
# including a few coins that changed rank position (rows 42/44/45) and a
# replacement coin in row 51 (ApeXProtocol -> Aave).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "46.912.45"
$ws.Range("E2").Value = "  +3.52%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "2.611.71"
$ws.Range("E3").Value = "  +6.99%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.30%  "

# Row 5: BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.26"
$ws.Range("E5").Value = "  +4.46%  "

# Row 6: Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.85"
$ws.Range("E6").Value = "  +7.01%  "

# Row 7: XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.601"
$ws.Range("E7").Value = "  +7.48%  "

# Row 8: USDC
$ws.Range("E8").Value = "  +0.07%  "

# Row 9: Cardano
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.584"
$ws.Range("E9").Value = "  +16.19%  "

# Row 10: Avalanche
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.50"
$ws.Range("E10").Value = "  +14.46%  "

# Row 11: OKB
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.95"
$ws.Range("E11").Value = "  +2.30%  "

# Row 12: Dogecoin
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0849"
$ws.Range("E12").Value = "  +8.92%  "

# Row 13: Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.20"
$ws.Range("E13").Value = "  +16.70%  "

# Row 14: WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "3.018.52"
$ws.Range("E14").Value = "  +7.39%  "

# Row 15: TRON
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.106"
$ws.Range("E15").Value = "  +1.82%  "

# Row 16: WrappedEther
$ws.Range("D16").Value = "2.621.75"
$ws.Range("E16").Value = "  +7.86%  "

# Row 17: Polygon
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.921"
$ws.Range("E17").Value = "  +9.18%  "

# Row 18: Chainlink
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "15.03"
$ws.Range("E18").Value = "  +5.87%  "

# Row 19: WrappedBTC
$ws.Range("D19").Value = "47.112.68"
$ws.Range("E19").Value = "  +4.08%  "

# Row 20: ShibaInu
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000101"
$ws.Range("E20").Value = "  +7.75%  "

# Row 21: InternetComputer(DFINITY)
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.95"
$ws.Range("E21").Value = "  +4.43%  "

# Row 22: Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.69"
$ws.Range("E22").Value = "  +7.60%  "

# Row 23: Litecoin
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.21"
$ws.Range("E23").Value = "  +7.93%  "

# Row 24: BitcoinCash
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "275.52"
$ws.Range("E24").Value = "  +15.05%  "

# Row 25: PancakeSwap
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.03"
$ws.Range("E25").Value = "  +9.13%  "

# Row 26: EthereumClassic
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "30.48"
$ws.Range("E26").Value = "  +41.89%  "

# Row 27: ImmutableX
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.17"
$ws.Range("E27").Value = "  +12.85%  "

# Row 28: Dai
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.997"
$ws.Range("E28").Value = "  -0.11%  "

# Row 29: LEO
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.04"
$ws.Range("E29").Value = "  +0.65%  "

# Row 30: Cosmos
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.54"
$ws.Range("E30").Value = "  +9.56%  "

# Row 31: Toncoin
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.31"
$ws.Range("E31").Value = "  +4.26%  "

# Row 32: InjectiveProtocol
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "38.87"
$ws.Range("E32").Value = "  +4.24%  "

# Row 33: Filecoin
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.11"
$ws.Range("E33").Value = "  +12.46%  "

# Row 34: LidoDAOToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.63"
$ws.Range("E34").Value = "  -6.74%  "

# Row 35: WEMIXToken
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.87"
$ws.Range("E35").Value = "  +5.39%  "

# Row 36: Hedera
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0843"
$ws.Range("E36").Value = "  +10.34%  "

# Row 37: Kaspa
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.126"
$ws.Range("E37").Value = "  +12.04%  "

# Row 38: ARBITRUM
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.19"
$ws.Range("E38").Value = "  +9.05%  "

# Row 39: Monero
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "150.53"
$ws.Range("E39").Value = "  +0.94%  "

# Row 40: Stellar
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.123"
$ws.Range("E40").Value = "  +6.75%  "

# Row 41: EnergySwap
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.52"
$ws.Range("E41").Value = "  +38.26%  "

# Row 42: NEARProtocol
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.66"
$ws.Range("E42").Value = "  +9.06%  "

# Row 43: RenderToken
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.14"
$ws.Range("E43").Value = "  +10.71%  "

# Row 44: Celestia
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0330"
$ws.Range("E44").Value = "  +12.29%  "

# Row 45: VeChain
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.63"
$ws.Range("E45").Value = "  +14.36%  "

# Row 46: Maker
$ws.Range("D46").Value = "2.166.27"
$ws.Range("E46").Value = "  +8.53%  "

# Row 47: BitcoinSV
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "96.13"
$ws.Range("E47").Value = "  +8.82%  "

# Row 48: FirstDigitalUSD
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  +0.10%  "

# Row 49: FraxShare
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.73"
$ws.Range("E49").Value = "  +13.15%  "

# Row 50: Stacks
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.78"
$ws.Range("E50").Value = "  +4.81%  "

# Row 51: ApeXProtocol
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "108.79"
$ws.Range("E51").Value = "  +6.64%  "

